$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$check = [char]0x2713

# Instruction names first (mov, mvn, orr) to seed the shared-string table in
# the same order the original author typed them, then their descriptions
# (note the description order differs: ORR, MOV, MVN), then andOp + its
# description last.
$ws.Range("A3").Value = "mov"
$ws.Range("A4").Value = "mvn"
$ws.Range("A5").Value = "orr"
$ws.Range("B5").Value = "Implements ORR function"
$ws.Range("B3").Value = "Implements MOV function"
$ws.Range("B4").Value = "Implements MVN function"
$ws.Range("A6").Value = "andOp"
$ws.Range("B6").Value = "Implements AND function"

# Row 3 - mov
$ws.Range("C3").Value = $check
$ws.Range("E3").Value = $check
$ws.Range("F3").Value = $check
$ws.Range("H3").Value = $check

# Row 4 - mvn
$ws.Range("C4").Value = $check
$ws.Range("E4").Value = $check
$ws.Range("F4").Value = $check
$ws.Range("H4").Value = $check

# Row 5 - orr
$ws.Range("C5").Value = $check
$ws.Range("D5").Value = $check
$ws.Range("E5").Value = $check
$ws.Range("F5").Value = $check
$ws.Range("H5").Value = $check

# Row 6 - andOp
$ws.Range("C6").Value = $check
$ws.Range("D6").Value = $check
$ws.Range("E6").Value = $check
$ws.Range("F6").Value = $check
$ws.Range("H6").Value = $check

$ws.Range("H6").Select()
